$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44461
$ws.Range("D6").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107002
$ws.Range("J6").Value = "Chirimoya"
$ws.Range("K6").Value = "Cultivar IV Región"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 29000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 29500
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 2950
$ws.Range("T6").Value = 10
